$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextCell $ws.Cells.Item(2, 4) "30.406.80"
Set-TextCell $ws.Cells.Item(2, 5) "  -2.00%  "

# Row 3
Set-TextCell $ws.Cells.Item(3, 4) "1.909.73"
Set-TextCell $ws.Cells.Item(3, 5) "  -2.54%  "

# Row 4
Set-TextCell $ws.Cells.Item(4, 4) "1.001"
Set-TextCell $ws.Cells.Item(4, 5) "  -0.01%  "

# Row 5
Set-TextCell $ws.Cells.Item(5, 4) "238.26"
Set-TextCell $ws.Cells.Item(5, 5) "  -2.58%  "

# Row 6
Set-TextCell $ws.Cells.Item(6, 4) "0.9997"
Set-TextCell $ws.Cells.Item(6, 5) "  -0.10%  "

# Row 7
Set-TextCell $ws.Cells.Item(7, 5) "  -2.33%  "

# Row 8
Set-TextCell $ws.Cells.Item(8, 4) "0.2855"
Set-TextCell $ws.Cells.Item(8, 5) "  -3.29%  "

# Row 9
Set-TextCell $ws.Cells.Item(9, 4) "0.06691"
Set-TextCell $ws.Cells.Item(9, 5) "  -4.20%  "

# Row 10
Set-TextCell $ws.Cells.Item(10, 4) "18.79"
Set-TextCell $ws.Cells.Item(10, 5) "  -3.47%  "

# Row 11
Set-TextCell $ws.Cells.Item(11, 4) "102.64"
Set-TextCell $ws.Cells.Item(11, 5) "  -4.89%  "

# Row 12
Set-TextCell $ws.Cells.Item(12, 4) "0.07717"
Set-TextCell $ws.Cells.Item(12, 5) "  -1.12%  "

# Row 13
Set-TextCell $ws.Cells.Item(13, 4) "1.908.42"
Set-TextCell $ws.Cells.Item(13, 5) "  -2.63%  "

# Row 14
Set-TextCell $ws.Cells.Item(14, 4) "5.196"
Set-TextCell $ws.Cells.Item(14, 5) "  -5.52%  "

# Row 15
Set-TextCell $ws.Cells.Item(15, 4) "0.6724"
Set-TextCell $ws.Cells.Item(15, 5) "  -4.12%  "

# Row 16
Set-TextCell $ws.Cells.Item(16, 4) "30.414.38"
Set-TextCell $ws.Cells.Item(16, 5) "  -2.03%  "

# Row 17
Set-TextCell $ws.Cells.Item(17, 4) "258.88"
Set-TextCell $ws.Cells.Item(17, 5) "  -7.88%  "

# Row 19
Set-TextCell $ws.Cells.Item(19, 4) "0.000007494"
Set-TextCell $ws.Cells.Item(19, 5) "  -3.73%  "

# Row 20
Set-TextCell $ws.Cells.Item(20, 4) "12.70"
Set-TextCell $ws.Cells.Item(20, 5) "  -4.72%  "

# Row 21
Set-TextCell $ws.Cells.Item(21, 4) "5.405"
Set-TextCell $ws.Cells.Item(21, 5) "  -2.49%  "

# Row 22
Set-TextCell $ws.Cells.Item(22, 4) "1.001"
Set-TextCell $ws.Cells.Item(22, 5) "  +0.09%  "

# Row 23
Set-TextCell $ws.Cells.Item(23, 2) "BitDAO"
Set-TextCell $ws.Cells.Item(23, 3) "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
Set-TextCell $ws.Cells.Item(23, 4) "0.4558"
Set-TextCell $ws.Cells.Item(23, 5) "  -8.37%  "

# Row 24
Set-TextCell $ws.Cells.Item(24, 2) "Chainlink"
Set-TextCell $ws.Cells.Item(24, 3) "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws.Cells.Item(24, 4) "6.292"
Set-TextCell $ws.Cells.Item(24, 5) "  -3.54%  "

# Row 25
Set-TextCell $ws.Cells.Item(25, 2) "Cosmos"
Set-TextCell $ws.Cells.Item(25, 3) "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws.Cells.Item(25, 4) "9.446"
Set-TextCell $ws.Cells.Item(25, 5) "  -4.24%  "

# Row 26
Set-TextCell $ws.Cells.Item(26, 2) "Monero"
Set-TextCell $ws.Cells.Item(26, 3) "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Cells.Item(26, 4) "164.47"
Set-TextCell $ws.Cells.Item(26, 5) "  -2.41%  "

# Row 27
Set-TextCell $ws.Cells.Item(27, 2) "EthereumClassic"
Set-TextCell $ws.Cells.Item(27, 3) "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Cells.Item(27, 4) "18.97"
Set-TextCell $ws.Cells.Item(27, 5) "  -4.99%  "

# Row 28
Set-TextCell $ws.Cells.Item(28, 2) "LidoDAOToken"
Set-TextCell $ws.Cells.Item(28, 3) "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws.Cells.Item(28, 4) "2.056"
Set-TextCell $ws.Cells.Item(28, 5) "  -6.30%  "

# Row 29
Set-TextCell $ws.Cells.Item(29, 2) "Stellar"
Set-TextCell $ws.Cells.Item(29, 3) "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Cells.Item(29, 4) "0.1009"
Set-TextCell $ws.Cells.Item(29, 5) "  -3.82%  "

# Row 30
Set-TextCell $ws.Cells.Item(30, 2) "Toncoin"
Set-TextCell $ws.Cells.Item(30, 3) "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Cells.Item(30, 4) "1.378"
Set-TextCell $ws.Cells.Item(30, 5) "  -0.74%  "

# Row 31
Set-TextCell $ws.Cells.Item(31, 2) "Filecoin"
Set-TextCell $ws.Cells.Item(31, 3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Cells.Item(31, 4) "4.626"
Set-TextCell $ws.Cells.Item(31, 5) "  -0.20%  "

# Row 32
Set-TextCell $ws.Cells.Item(32, 2) "PancakeSwap"
Set-TextCell $ws.Cells.Item(32, 3) "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Cells.Item(32, 4) "1.510"
Set-TextCell $ws.Cells.Item(32, 5) "  -3.76%  "

# Row 33
Set-TextCell $ws.Cells.Item(33, 2) "InternetComputer(DFINITY)"
Set-TextCell $ws.Cells.Item(33, 3) "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Cells.Item(33, 4) "4.253"
Set-TextCell $ws.Cells.Item(33, 5) "  -4.85%  "

# Row 34
Set-TextCell $ws.Cells.Item(34, 2) "Hedera"
Set-TextCell $ws.Cells.Item(34, 3) "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Cells.Item(34, 4) "0.04787"
Set-TextCell $ws.Cells.Item(34, 5) "  -2.78%  "

# Row 35
Set-TextCell $ws.Cells.Item(35, 2) "ImmutableX"
Set-TextCell $ws.Cells.Item(35, 3) "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws.Cells.Item(35, 4) "0.7307"
Set-TextCell $ws.Cells.Item(35, 5) "  -3.35%  "

# Row 36
Set-TextCell $ws.Cells.Item(36, 2) "ARBITRUM"
Set-TextCell $ws.Cells.Item(36, 3) "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws.Cells.Item(36, 4) "1.111"
Set-TextCell $ws.Cells.Item(36, 5) "  -5.18%  "

# Row 37
Set-TextCell $ws.Cells.Item(37, 2) "Frax"
Set-TextCell $ws.Cells.Item(37, 3) "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws.Cells.Item(37, 4) "0.9992"
Set-TextCell $ws.Cells.Item(37, 5) "  -0.15%  "

# Row 38
Set-TextCell $ws.Cells.Item(38, 2) "HuobiToken"
Set-TextCell $ws.Cells.Item(38, 3) "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws.Cells.Item(38, 4) "2.709"
Set-TextCell $ws.Cells.Item(38, 5) "  -0.91%  "

# Row 39
Set-TextCell $ws.Cells.Item(39, 2) "VeChain"
Set-TextCell $ws.Cells.Item(39, 3) "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Cells.Item(39, 4) "0.01927"
Set-TextCell $ws.Cells.Item(39, 5) "  -4.19%  "

# Row 40
Set-TextCell $ws.Cells.Item(40, 2) "MXToken"
Set-TextCell $ws.Cells.Item(40, 3) "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws.Cells.Item(40, 4) "2.602"
Set-TextCell $ws.Cells.Item(40, 5) "  -3.77%  "

# Row 41
Set-TextCell $ws.Cells.Item(41, 2) "FraxShare"
Set-TextCell $ws.Cells.Item(41, 3) "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws.Cells.Item(41, 4) "6.230"
Set-TextCell $ws.Cells.Item(41, 5) "  -5.13%  "

# Row 42
Set-TextCell $ws.Cells.Item(42, 2) "Aave"
Set-TextCell $ws.Cells.Item(42, 3) "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Cells.Item(42, 4) "74.71"
Set-TextCell $ws.Cells.Item(42, 5) "  -4.48%  "

# Row 43
Set-TextCell $ws.Cells.Item(43, 2) "RenderToken"
Set-TextCell $ws.Cells.Item(43, 3) "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Cells.Item(43, 4) "1.989"
Set-TextCell $ws.Cells.Item(43, 5) "  -6.96%  "

# Row 44
Set-TextCell $ws.Cells.Item(44, 2) "TrustWalletToken"
Set-TextCell $ws.Cells.Item(44, 3) "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Cells.Item(44, 4) "0.8649"
Set-TextCell $ws.Cells.Item(44, 5) "  -4.18%  "

# Row 45
Set-TextCell $ws.Cells.Item(45, 2) "Quant"
Set-TextCell $ws.Cells.Item(45, 3) "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws.Cells.Item(45, 4) "106.67"
Set-TextCell $ws.Cells.Item(45, 5) "  -2.55%  "

# Row 46
Set-TextCell $ws.Cells.Item(46, 2) "Maker"
Set-TextCell $ws.Cells.Item(46, 3) "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Cells.Item(46, 4) "1.065.09"
Set-TextCell $ws.Cells.Item(46, 5) "  +5.33%  "

# Row 47
Set-TextCell $ws.Cells.Item(47, 2) "TheSandbox"
Set-TextCell $ws.Cells.Item(47, 3) "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws.Cells.Item(47, 4) "0.4248"
Set-TextCell $ws.Cells.Item(47, 5) "  -4.65%  "

# Row 48
Set-TextCell $ws.Cells.Item(48, 2) "PaxDollar"
Set-TextCell $ws.Cells.Item(48, 3) "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws.Cells.Item(48, 4) "0.9989"
Set-TextCell $ws.Cells.Item(48, 5) "  -0.22%  "

# Row 49
Set-TextCell $ws.Cells.Item(49, 2) "Aptos"
Set-TextCell $ws.Cells.Item(49, 3) "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Cells.Item(49, 4) "7.478"
Set-TextCell $ws.Cells.Item(49, 5) "  -8.13%  "

# Row 50
Set-TextCell $ws.Cells.Item(50, 2) "Algorand"
Set-TextCell $ws.Cells.Item(50, 3) "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws.Cells.Item(50, 4) "0.1198"
Set-TextCell $ws.Cells.Item(50, 5) "  -4.64%  "

# Row 51
Set-TextCell $ws.Cells.Item(51, 2) "Elrond"
Set-TextCell $ws.Cells.Item(51, 3) "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell $ws.Cells.Item(51, 4) "35.05"
Set-TextCell $ws.Cells.Item(51, 5) "  -2.72%  "
